$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update coin price/volume data to reflect the latest scrape ---

$style = $ws.Range("D2").Style
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "39.913.93"
$ws.Range("D2").Style = $style
$ws.Range("E2").Value = "  +0.02%  "
$style = $ws.Range("D3").Style
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.205.84"
$ws.Range("D3").Style = $style
$ws.Range("E3").Value = "  -0.80%  "
$style = $ws.Range("D4").Style
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = $style
$ws.Range("E4").Value = "  -0.07%  "
$style = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "293.08"
$ws.Range("D5").Style = $style
$ws.Range("E5").Value = "  -0.51%  "
$style = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "87.14"
$ws.Range("D6").Style = $style
$ws.Range("E6").Value = "  +1.74%  "
$style = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.508"
$ws.Range("D7").Style = $style
$ws.Range("E7").Value = "  -0.96%  "
$ws.Range("E8").Value = "  -0.01%  "
$style = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.471"
$ws.Range("D9").Style = $style
$ws.Range("E9").Value = "  +0.19%  "
$style = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "29.85"
$ws.Range("D10").Style = $style
$ws.Range("E10").Value = "  -3.66%  "
$style = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0775"
$ws.Range("D11").Style = $style
$ws.Range("E11").Value = "  -1.65%  "
$style = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "49.49"
$ws.Range("D12").Style = $style
$ws.Range("E12").Value = "  +5.59%  "
$ws.Range("E13").Value = "  +2.60%  "
$style = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.43"
$ws.Range("D14").Style = $style
$ws.Range("E14").Value = "  +0.19%  "
$style = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.549.01"
$ws.Range("D15").Style = $style
$ws.Range("E15").Value = "  -0.92%  "
$style = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "13.71"
$ws.Range("D16").Style = $style
$ws.Range("E16").Value = "  -2.22%  "
$style = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.199.18"
$ws.Range("D17").Style = $style
$ws.Range("E17").Value = "  -1.18%  "
$style = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.724"
$ws.Range("D18").Style = $style
$ws.Range("E18").Value = "  -0.35%  "
$style = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "39.842.86"
$ws.Range("D19").Style = $style
$ws.Range("E19").Value = "  -0.04%  "
$style = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0883"
$ws.Range("D20").Style = $style
$ws.Range("E20").Value = "  -0.42%  "
$style = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.27"
$ws.Range("D21").Style = $style
$ws.Range("E21").Value = "  +5.14%  "
$style = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.74"
$ws.Range("D22").Style = $style
$ws.Range("E22").Value = "  -0.68%  "
$style = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "65.18"
$ws.Range("D23").Style = $style
$ws.Range("E23").Value = "  +0.15%  "
$style = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "236.45"
$ws.Range("D24").Style = $style
$ws.Range("E24").Value = "  +0.67%  "
$ws.Range("E25").Value = "  +0.10%  "
$style = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.45"
$ws.Range("D26").Style = $style
$ws.Range("E26").Value = "  -0.22%  "
$ws.Range("E27").Value = "  -2.03%  "
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$style = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "22.42"
$ws.Range("D28").Style = $style
$ws.Range("E28").Value = "  -1.21%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$style = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.15"
$ws.Range("D29").Style = $style
$ws.Range("E29").Value = "  -3.38%  "
$style = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.15"
$ws.Range("D30").Style = $style
$ws.Range("E30").Value = "  -0.36%  "
$style = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "155.25"
$ws.Range("D31").Style = $style
$ws.Range("E31").Value = "  +2.36%  "
$style = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "31.53"
$ws.Range("D32").Style = $style
$ws.Range("E32").Value = "  -4.80%  "
$ws.Range("E33").Value = "  -0.11%  "
$style = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.88"
$ws.Range("D34").Style = $style
$ws.Range("E34").Value = "  +0.45%  "
$style = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0708"
$ws.Range("D35").Style = $style
$ws.Range("E35").Value = "  -0.78%  "
$ws.Range("E36").Value = "  -2.34%  "
$style = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.83"
$ws.Range("D37").Style = $style
$ws.Range("E37").Value = "  +4.55%  "
$style = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.112"
$ws.Range("D38").Style = $style
$ws.Range("E38").Value = "  +0.47%  "
$style = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0976"
$ws.Range("D39").Style = $style
$ws.Range("E39").Value = "  -1.81%  "
$style = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "15.36"
$ws.Range("D40").Style = $style
$ws.Range("E40").Value = "  -5.37%  "
$ws.Range("E41").Value = "  -0.92%  "
$style = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.121.81"
$ws.Range("D42").Style = $style
$ws.Range("E42").Value = "  +4.73%  "
$style = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.73"
$ws.Range("D43").Style = $style
$ws.Range("E43").Value = "  -1.75%  "
$ws.Range("E44").Value = "  -3.38%  "
$ws.Range("E45").Value = "  -0.91%  "
$style = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "17.56"
$ws.Range("D46").Style = $style
$ws.Range("E46").Value = "  +7.66%  "
$style = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.63"
$ws.Range("D47").Style = $style
$ws.Range("E47").Value = "  -3.47%  "
$style = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.64"
$ws.Range("D48").Style = $style
$ws.Range("E48").Value = "  +3.06%  "
$style = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.416.30"
$ws.Range("D49").Style = $style
$ws.Range("E49").Value = "  -1.12%  "
$ws.Range("E50").Value = "  +1.26%  "
$style = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.10"
$ws.Range("D51").Style = $style
$ws.Range("E51").Value = "  +0.95%  "
